$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.31 = 12698.68 pesos`n✅ 12698.68 pesos = 3.29 = 937.65 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 302
$wsTasas.Range("O10").Value = 3835
$wsTasas.Range("N12").Value = 3859.8
$wsTasas.Range("O12").Value = 285
